# Csv-to-csv transcoding; suppress LO messages
#
# Sheet1: append a new data row (Françoise-Athénaïs / de Rochechouart /
#         -94683), center the Date column's values, resize a few columns,
#         and make Sheet1 the active/selected sheet.
# "Another Sheet": its shared-string-backed values don't change in
#         meaning, just stop being the active tab and pick up the
#         new selection left behind on Sheet1's column D.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Another Sheet")

# --- new row of data on Sheet1 -------------------------------------------
$ws1.Range("A4").Value = "Françoise-Athénaïs"
$ws1.Range("B4").Value = "de Rochechouart"
$ws1.Range("D4").Value = -94683

# Center-align the Date column's values (D2:D4), same number format.
$ws1.Range("D2:D4").HorizontalAlignment = -4108

# --- column widths on Sheet1 ----------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 16.42
$ws1.Columns.Item(2).ColumnWidth = 15.09
$ws1.Columns.Item(4).ColumnWidth = 11.42

# --- selections / active sheet --------------------------------------------
$ws2.Range("G5").Select()
$ws2.Columns.Item(4).Select()
$ws1.Range("D4").Select()
$ws1.Columns.Item(4).Select()
$ws1.Activate()
